# Insert a new "Lab 400: Querying External Data" slide between the
# existing "Lab 300" slide and the existing "Lab 500" slide.
#
# Implementation: duplicate the current "Lab 500" slide (slide 5). The
# duplicate lands right after the original (new slide 6). We then:
#   - turn the ORIGINAL (still at position 5) into the new "Lab 400" slide,
#     resizing/repositioning its text box and editing its text in place
#     (mirrors how "Lab 300" etc. were authored: "Lab " + number + ":" as
#     separate runs, then a second paragraph with the lab title).
#   - reset the DUPLICATE's (position 6) text box back to the standard
#     textbox geometry shared by the other lab slides, keeping the
#     original "Lab 500: Managing and Monitoring via Cloud Console" text.

$p = $ppt.ActivePresentation

# "Lab 500" slide (Test Drive / Workshop / Lab 500: ... ) is slide 5.
$labSlide = $p.Slides.Item(5)

# Duplicate it -> new slide inserted right after (position 6).
$dup = $labSlide.Duplicate()

# --- Slide 5 (original) becomes "Lab 400: Querying External Data" ---
$shp400 = $labSlide.Shapes.Item(3)
$tr400 = $shp400.TextFrame.TextRange

$full = $tr400.Text
$idx = $full.IndexOf("500")
$tr400.Characters($idx + 1, 3).Text = "400"

$full = $tr400.Text
$oldTitle = "Managing and Monitoring via Cloud Console"
$idx2 = $full.IndexOf($oldTitle)
$tr400.Characters($idx2 + 1, $oldTitle.Length).Text = "Querying External Data"

# Resize / reposition the text box to its new custom geometry.
$shp400.Left = 499.2896062992126
$shp400.Top = 101.221968503937
$shp400.Width = 450.0437007874016
$shp400.Height = 273.84842519685037

# --- Slide 6 (duplicate) stays "Lab 500" but reverts to the standard box ---
$s500 = $p.Slides.Item(6)
$shp500 = $s500.Shapes.Item(3)
$shp500.Left = 499.28968503937006
$shp500.Top = 101.221968503937
$shp500.Width = 352.59811023622046
$shp500.Height = 327.164094488189
